# Add new "Table S2 - PERMANOVA AIC" sheet (PCA/PERMANOVA AIC comparison
# table), shifting the existing Table S2-S8 sheets down to S3-S9, and
# refresh a handful of PERMANOVA p-values that were recalculated as part
# of the same analysis update.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the existing S2..S8 sheets up to S3..S9 first (process from
#    the end backwards so no two sheets ever collide on a name).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Table S8 - HostVsymb Plast GLM").Name  = "Table S9 - HostVsymb Plast GLM"
$wb.Worksheets.Item("Table S7 - HostVsymb Plast AIC").Name  = "Table S8 - HostVsymb Plast AIC"
$wb.Worksheets.Item("Table S6 - HostVsymb PERMANOVA").Name  = "Table S7 - HostVsymb PERMANOVA"
$wb.Worksheets.Item("Table S5 - Species PERMANOVA").Name    = "Table S6 - Species PERMANOVA"
$wb.Worksheets.Item("Table S4 - Plasticity GLM").Name       = "Table S5 - Plasticity GLM"
$wb.Worksheets.Item("Table S3 - PERMANOVA").Name             = "Table S4 - PERMANOVA"
$wb.Worksheets.Item("Table S2 - Plasticity AIC").Name        = "Table S3 - Plasticity AIC"

# ---------------------------------------------------------------------
# 2) Insert the brand-new sheet right after "Table S1 - Sample Size" and
#    populate it with the PERMANOVA AIC comparison table.
# ---------------------------------------------------------------------
$after = $wb.Worksheets.Item("Table S1 - Sample Size")
$newSheet = $wb.Worksheets.Add($null, $after)
$newSheet.Name = "Table S2 - PERMANOVA AIC"

$newSheet.Cells.Item(1,1).Value = "group"
$newSheet.Cells.Item(1,2).Value = "Full interactive model AIC"
$newSheet.Cells.Item(1,3).Value = "Best fit (additive) model AIC"

$newSheet.Cells.Item(2,1).Value = "S. siderea"
$newSheet.Cells.Item(2,2).Value = 692.3
$newSheet.Cells.Item(2,3).Value = 678.5

$newSheet.Cells.Item(3,1).Value = "P. strigosa"
$newSheet.Cells.Item(3,2).Value = 696.9
$newSheet.Cells.Item(3,3).Value = 685.8

$newSheet.Cells.Item(4,1).Value = "P. astreoides"
$newSheet.Cells.Item(4,2).Value = 500.5
$newSheet.Cells.Item(4,3).Value = 497.1

$newSheet.Columns.Item(1).ColumnWidth = 13.71
$newSheet.Columns.Item(2).ColumnWidth = 26.71
$newSheet.Columns.Item(3).ColumnWidth = 29.71

# ---------------------------------------------------------------------
# 3) Refresh the handful of recalculated PERMANOVA p-values in the
#    tables that are now S4, S6 and S7.
# ---------------------------------------------------------------------

# Table S4 - PERMANOVA
$s4 = $wb.Worksheets.Item("Table S4 - PERMANOVA")
$s4.Range("F3").Value = 0.05463
$s4.Range("F9").Value = 0.22252
$s4.Range("F12").Value = 0.48568

# Table S6 - Species PERMANOVA
$s6 = $wb.Worksheets.Item("Table S6 - Species PERMANOVA")
$s6.Range("F3").Value = 0.0966
$s6.Range("F4").Value = 0.00133
$s6.Range("F7").Value = 0.02065
$s6.Range("F8").Value = 0.00266

# Table S7 - HostVsymb PERMANOVA
$s7 = $wb.Worksheets.Item("Table S7 - HostVsymb PERMANOVA")
$s7.Range("F2").Value  = 0.77615
$s7.Range("F3").Value  = 0.002
$s7.Range("K3").Value  = 0.09061
$s7.Range("F4").Value  = 0.58894
$s7.Range("F7").Value  = 0.26382
$s7.Range("K7").Value  = 0.27848
$s7.Range("F9").Value  = 0.17255
$s7.Range("K9").Value  = 0.002
$s7.Range("F12").Value = 0.01666
$s7.Range("F13").Value = 0.10193
$s7.Range("F14").Value = 0.1992
$s7.Range("K14").Value = 0.48168
